# Updating barcode_offset and umi_offset
#
# 1. barcode_offset list sheet gains a third offset option "10,48,86" and
#    the original "0" option becomes "0,38,76".
# 2. umi_offset list sheet gains a new first option "1" (existing options
#    shift down).
# 3. The RNAseq sheet's data-validation list ranges for the barcode_offset
#    and umi_offset columns are widened to include the new rows.
# 4. The .metadata sheet's pav:createdOn timestamp is bumped.

$wb = $excel.ActiveWorkbook

# --- 1. barcode_offset ------------------------------------------------
$boSheet = $wb.Worksheets.Item("barcode_offset")
$boSheet.Range("A1").Value = "0,38,76"
$boSheet.Range("A5").Value = "10,48,86"

# --- 2. umi_offset ------------------------------------------------------
$uoSheet = $wb.Worksheets.Item("umi_offset")
$uoSheet.Range("A1").Insert()
# Leading apostrophe forces the numeric-looking "1" to be stored as text
# (matching the other cells in this list, which are all text); ClearFormats
# drops the "quote prefix" style flag that the apostrophe entry leaves
# behind, so no stray style gets attached to the cell.
$uoSheet.Range("A1").Value = "'1"
$uoSheet.Range("A1").ClearFormats()

# --- 3. RNAseq data validations ------------------------------------------
$main = $wb.Worksheets.Item("RNAseq")

$dvBarcode = $main.Range("O2:O1001").Validation
$dvBarcode.Modify(3, 1, 1, "='barcode_offset'!`$A`$1:`$A`$5")
$dvBarcode.IgnoreBlank = $true
$dvBarcode.ShowError = $true
$dvBarcode.ShowInput = $false
$dvBarcode.ErrorTitle = "Validation Error"
$dvBarcode.ErrorMessage = ""

$dvUmi = $main.Range("R2:R1001").Validation
$dvUmi.Modify(3, 1, 1, "='umi_offset'!`$A`$1:`$A`$3")
$dvUmi.IgnoreBlank = $true
$dvUmi.ShowError = $true
$dvUmi.ShowInput = $false
$dvUmi.ErrorTitle = "Validation Error"
$dvUmi.ErrorMessage = ""

# --- 4. .metadata pav:createdOn -----------------------------------------
$meta = $wb.Worksheets.Item(".metadata")
$meta.Cells.Item(2, 3).Value = "2023-10-31T13:53:33-07:00"
